$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New armor / equipment rows (rows 8-24), columns:
#   B = ARMOR name, C = TYPE, D = CLASS (all new items are "Exotic" rarity-class)
$ws.Range("B8").Value = "Psijic Warden Helm"
$ws.Range("C8").Value = "Heavy Helm"
$ws.Range("D8").Value = "Exotic"

$ws.Range("B9").Value = "Psijic Warden Cuirass"
$ws.Range("C9").Value = "Heavy Cuirass"
$ws.Range("D9").Value = "Exotic"

$ws.Range("B10").Value = "Psijic Warden Gauntlets"
$ws.Range("C10").Value = "Heavy Gauntlets"
$ws.Range("D10").Value = "Exotic"

$ws.Range("B11").Value = "Psijic Warden Boots"
$ws.Range("C11").Value = "Heavy Boots"
$ws.Range("D11").Value = "Exotic"

$ws.Range("B12").Value = "Psijic Warden Aegis"
$ws.Range("C12").Value = "Heavy Shield"
$ws.Range("D12").Value = "Exotic"

$ws.Range("B13").Value = "Psijic Agent Robe"
$ws.Range("C13").Value = "Clothing"
$ws.Range("D13").Value = "Exotic"

$ws.Range("B14").Value = "Psijic Agent Hood"
$ws.Range("C14").Value = "Clothing"
$ws.Range("D14").Value = "Exotic"

$ws.Range("B15").Value = "Psijic Ritemaster Headdress"
$ws.Range("C15").Value = "Clothing"
$ws.Range("D15").Value = "Exotic"

$ws.Range("B16").Value = "Psijic Ritemaster Regalia"
$ws.Range("C16").Value = "Clothing"
$ws.Range("D16").Value = "Exotic"

$ws.Range("B17").Value = "Psijic Mystical Raiment"
$ws.Range("C17").Value = "Clothing"
$ws.Range("D17").Value = "Exotic"

$ws.Range("B18").Value = "Psijic Apprentice Robe"
$ws.Range("C18").Value = "Clothing"
$ws.Range("D18").Value = "Exotic"

$ws.Range("B19").Value = "Psijic Apprentice Hood"
$ws.Range("C19").Value = "Clothing"
$ws.Range("D19").Value = "Exotic"

$ws.Range("B20").Value = "Psijic Sage Robe"
$ws.Range("C20").Value = "Clothing"
$ws.Range("D20").Value = "Exotic"

$ws.Range("B21").Value = "Psijic Sage Hood"
$ws.Range("C21").Value = "Clothing"
$ws.Range("D21").Value = "Exotic"

$ws.Range("B22").Value = "Psijic Counselor Regalia"
$ws.Range("C22").Value = "Clothing"
$ws.Range("D22").Value = "Exotic"

$ws.Range("B23").Value = "Psijic Counselor Headdress"
$ws.Range("C23").Value = "Clothing"
$ws.Range("D23").Value = "Exotic"

$ws.Range("B24").Value = "Psijic Keeper Robe"
$ws.Range("C24").Value = "Clothing"
$ws.Range("D24").Value = "Exotic"

# Widen the ARMOR and TYPE columns to fit their new, longer contents.
$ws.Columns.Item(2).ColumnWidth = 25.14
$ws.Columns.Item(3).ColumnWidth = 14.6

# Leave the selection where the author left off.
$ws.Range("E21").Select()
